$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B82").Value = "Pop-up on device overview screen if device is ready for setup"
$ws.Range("C82").Value = "Open"

$ws.Range("B83").Value = "Intercom should continue to work even when cloud connection is lost"
$ws.Range("C83").Value = "Open"

$ws.Range("B84").Value = "Remove blank screen from app"
$ws.Range("C84").Value = "Open"

$ws.Range("B85").Select()
